$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = -21.107
$ws.Range("A10").Value = -20.945
$ws.Range("A12").Value = -21.452
$ws.Range("E13").Value = 12.817
$ws.Range("A18").Value = -21.766
$ws.Range("A25").Value = -21.534
